# Test12.xlsx edit:
#  - On "Not Normalized", insert 3 blank rows above the existing table
#    (the table, originally A1:E9, becomes A4:E12) and move the selection
#    to J11.
#  - Append four new, empty worksheets at the end of the workbook:
#    "New Worksheet_914000", "New Worksheet_345000", "New Worksheet_928000",
#    "New Worksheet_256000".
#  - Keep "Not Normalized" the active/selected sheet (adding sheets makes
#    the newest one active, so it is reactivated at the end).

$wb = $excel.ActiveWorkbook

$wsNotNormalized = $wb.Worksheets.Item("Not Normalized")
$null = $wsNotNormalized.Rows("1:3").Insert()

$newSheetNames = @(
    "New Worksheet_914000",
    "New Worksheet_345000",
    "New Worksheet_928000",
    "New Worksheet_256000"
)

foreach ($name in $newSheetNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $newSheet.Name = $name
}

$wsNotNormalized.Activate() | Out-Null
$wsNotNormalized.Range("J11").Select() | Out-Null

Write-Output "Inserted 3 rows on 'Not Normalized' and added sheets: $($newSheetNames -join ', ')"
